$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ASML)
$ws.Range("D2").Value = 1107.73
$ws.Range("E2").Value = 64.59999999999999
$ws.Range("F2").Value = 4.5
$ws.Range("N2").Value = 50.68470204858703

# Row 3 (TSM)
$ws.Range("D3").Value = 296.26
$ws.Range("E3").Value = 60.9
$ws.Range("F3").Value = 1.63
$ws.Range("N3").Value = 50.68470204858703

# Row 4 (AMD)
$ws.Range("D4").Value = 217.76
$ws.Range("E4").Value = 33.3
$ws.Range("F4").Value = 0.11
$ws.Range("N4").Value = 50.68470204858703

# Row 5 (NVDA)
$ws.Range("D5").Value = 181.91
$ws.Range("E5").Value = 41.3
$ws.Range("F5").Value = 2.78
$ws.Range("N5").Value = 50.68470204858703

# Row 6 (QCOM)
$ws.Range("D6").Value = 175.37
$ws.Range("E6").Value = 53.1
$ws.Range("F6").Value = 4.86
$ws.Range("I6").Value = 36
$ws.Range("K6").Value = 47.6
$ws.Range("N6").Value = 50.68470204858703

$wb.Save()
